$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Markdown table gained four new rows. Insert blank rows at the
# positions (in ascending order) where the new entries belong so the
# existing rows are pushed down exactly as in the diff.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(19).Insert()

# Populate the newly inserted rows with the new entries.
$ws.Range("A13").Value = '''2005'
$ws.Range("B13").Value = '**広島市衛生研究所** <br> [鶏肉からのカンピロバクターの定量および定性検査法の有効性評価](https://www.city.hiroshima.lg.jp/_res/projects/default_project/_page_/001/023/132/57088.pdf) <br>（広島市衛生研究所年報, 25, pp.44-46, 2006）'
$ws.Range("C13").Value = '未登録'

$ws.Range("A14").Value = '''2005'
$ws.Range("B14").Value = '**富山県衛生研究所** <br> [富山県におけるカンピロバクター分離状況(2005年)](https://www.pref.toyama.jp/documents/13568/nenpou29h17_1.pdf) <br>（富山県衛生研究所年報, 29 (2006), pp.174-177）'
$ws.Range("C14").Value = '未登録'

$ws.Range("A16").Value = '2004-2005'
$ws.Range("B16").Value = '**宮城県保健環境センター** <br> [鶏肉からの効率的なカンピロバクターの分離の検討と分離菌の性状](https://www.pref.miyagi.jp/documents/1943/617297.pdf) <br>（宮城県保健環境センター年報, 第24号, pp.117~120,2006）'
$ws.Range("C16").Value = '未登録'

$ws.Range("A19").Value = '2002-2003'
$ws.Range("B19").Value = '**香川県環境保健研究センター** <br> [鶏肉における _Campylobacter_ および _Salmonella_ の汚染状況](https://www.pref.kagawa.lg.jp/documents/2480/s8cm2a170906165025_f23_1.pdf) <br>（香川県環境保健研究センター所報 第3号, 2004, pp.187-190）'
$ws.Range("C19").Value = '未登録'
